# create register event crud
# Rebuild the "Hoja 1" sheet from a 2-column (id/Nombre) test sheet into a
# 7-column events table (id_event, name, description, date, hour, address,
# id_city) with 3 sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stamp the existing cell format (style index 1 - Arial/theme font) onto
#    the whole A1:G4 block first, via a formats-only paste, so every new
#    cell we are about to populate shares the same style as the original
#    A1/B1/A2/B2 cells instead of falling back to the default style.
$ws.Range("A1").Copy()
$ws.Range("A1:G4").PasteSpecial(-4122)  # xlPasteFormats

# 2) Header row
$ws.Range("A1").Value = "id_event"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "hour"
$ws.Range("F1").Value = "address"
$ws.Range("G1").Value = "id_city"

# 3) Data rows - the "date" column is forced to text (NumberFormat "@")
#    before the value is written so the date-like strings aren't coerced
#    into date serial numbers.
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Prueba edited1"
$ws.Range("C2").Value = "Prueba load 1 description"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2024-11-12"
$ws.Range("E2").Value = "3:00 pm"
$ws.Range("F2").Value = "Calle 1"
$ws.Range("G2").Value = 1

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Prueba edited 2"
$ws.Range("C3").Value = "Prueba load 2 descriptio"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2024-10-13"
$ws.Range("F3").Value = "Calle 2"
$ws.Range("G3").Value = 1

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Prueba edited 3"
$ws.Range("C4").Value = "Prueba load 3 description"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2024-12-14"
$ws.Range("E4").Value = "3:00 pm"
$ws.Range("F4").Value = "Calle 3"
$ws.Range("G4").Value = 1

# 4) Widen column B ("name") - 21.64 chars round-trips through Excel's
#    pixel-snapping to the OOXML width="22.5" recorded in the target file.
$ws.Columns.Item(2).ColumnWidth = 21.64
